$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DashboardPageData")
$ws.Activate()

# ---------------------------------------------------------------------------
# 1. Swap the runMode values of the two existing "testQuotesDashboardUI" rows
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = "Y"
$ws.Range("A4").Value = "N"

# ---------------------------------------------------------------------------
# 2. Append the three new test-data blocks starting at row 32. Formats are
#    copied from the matching existing section headers / data rows so the
#    new cells end up sharing the same cell styles already used elsewhere
#    on the sheet (yellow header band, bordered data cells, text-formatted
#    date cells).
# ---------------------------------------------------------------------------

# --- block: testBrokerFilteringSubmissionsList (rows 32-35) ---------------
$ws.Range("A1:B1").Copy()
$ws.Range("A32:B32").PasteSpecial(-4122)
$ws.Range("A32").Value = "testBrokerFilteringSubmissionsList"

$ws.Range("A2:D2").Copy()
$ws.Range("A33:D33").PasteSpecial(-4122)
$ws.Range("A33:D33").Copy()
$ws.Range("E33:I33").PasteSpecial(-4122)

$ws.Range("A2:D2").Copy()
$ws.Range("A34:D34").PasteSpecial(-4122)
$ws.Range("A34:D34").Copy()
$ws.Range("E34:H34").PasteSpecial(-4122)
$ws.Range("E29").Copy()
$ws.Range("I34").PasteSpecial(-4122)

$ws.Range("A2:D2").Copy()
$ws.Range("A35:D35").PasteSpecial(-4122)
$ws.Range("A35:D35").Copy()
$ws.Range("E35:H35").PasteSpecial(-4122)
$ws.Range("E29").Copy()
$ws.Range("I35").PasteSpecial(-4122)

$ws.Range("A33").Value = "runMode"
$ws.Range("B33").Value = "brokerId"
$ws.Range("C33").Value = "agentId"
$ws.Range("D33").Value = "agencyOfficeId"
$ws.Range("E33").Value = "productName"
$ws.Range("F33").Value = "status"
$ws.Range("G33").Value = "allProducts"
$ws.Range("H33").Value = "allStatuses"
$ws.Range("I33").Value = "endDate"

$ws.Range("A34").Value = "Y"
$ws.Range("B34").Value = 20217
$ws.Range("C34").Value = 237
$ws.Range("D34").Value = 8006
$ws.Range("E34").Value = "NetGuard® Plus"
$ws.Range("F34").Value = "Active"
$ws.Range("G34").Value = "All Products"
$ws.Range("H34").Value = "All Statuses"
$ws.Range("I34").Value = "12/31/2021"

$ws.Range("A35").Value = "N"
$ws.Range("B35").Value = 25997
$ws.Range("C35").Value = 7166
$ws.Range("D35").Value = 8414
$ws.Range("E35").Value = "NetGuard® Plus;QA Program 5204"
$ws.Range("F35").Value = "Active;Cancelled;Declined"
$ws.Range("G35").Value = "All Products"
$ws.Range("H35").Value = "All Statuses"
$ws.Range("I35").Value = "12/31/2021"

# --- block: testBrokerFilteringPoliciesList (rows 38-41) -------------------
$ws.Range("A1:B1").Copy()
$ws.Range("A38:B38").PasteSpecial(-4122)
$ws.Range("A38").Value = "testBrokerFilteringPoliciesList"

$ws.Range("A2:D2").Copy()
$ws.Range("A39:D39").PasteSpecial(-4122)
$ws.Range("A39:D39").Copy()
$ws.Range("E39:H39").PasteSpecial(-4122)

$ws.Range("A2:D2").Copy()
$ws.Range("A40:D40").PasteSpecial(-4122)
$ws.Range("A40:D40").Copy()
$ws.Range("E40:G40").PasteSpecial(-4122)
$ws.Range("E29").Copy()
$ws.Range("H40").PasteSpecial(-4122)

$ws.Range("A2:D2").Copy()
$ws.Range("A41:D41").PasteSpecial(-4122)
$ws.Range("A41:D41").Copy()
$ws.Range("E41:G41").PasteSpecial(-4122)
$ws.Range("E29").Copy()
$ws.Range("H41").PasteSpecial(-4122)

$ws.Range("A39").Value = "runMode"
$ws.Range("B39").Value = "brokerId"
$ws.Range("C39").Value = "agentId"
$ws.Range("D39").Value = "agencyOfficeId"
$ws.Range("E39").Value = "productName"
$ws.Range("F39").Value = "status"
$ws.Range("G39").Value = "allStatuses"
$ws.Range("H39").Value = "endDate"

$ws.Range("A40").Value = "Y"
$ws.Range("B40").Value = 20217
$ws.Range("C40").Value = 237
$ws.Range("D40").Value = 8006
$ws.Range("E40").Value = "NetGuard® Plus;Generic Program - NetGuard® Plus"
$ws.Range("F40").Value = "Renewed;Renewal Started;Expired"
$ws.Range("G40").Value = "All Statuses"
$ws.Range("H40").Value = "12/30/2021"

$ws.Range("A41").Value = "N"
$ws.Range("B41").Value = 25997
$ws.Range("C41").Value = 7166
$ws.Range("D41").Value = 8414
$ws.Range("E41").Value = "NetGuard® Plus;Generic Program - NetGuard® Plus"
$ws.Range("F41").Value = "Renewed;Renewal Started;Expired"
$ws.Range("G41").Value = "All Statuses"
$ws.Range("H41").Value = "12/30/2021"

# --- block: testPresenceOfContinueButtonOnQuotes (rows 44-47) --------------
$ws.Range("A1:B1").Copy()
$ws.Range("A44:B44").PasteSpecial(-4122)
$ws.Range("A44").Value = "testPresenceOfContinueButtonOnQuotes"

$ws.Range("A2:D2").Copy()
$ws.Range("A45:D45").PasteSpecial(-4122)
$ws.Range("A2:D2").Copy()
$ws.Range("A46:D46").PasteSpecial(-4122)
$ws.Range("A2:D2").Copy()
$ws.Range("A47:D47").PasteSpecial(-4122)

$ws.Range("A45").Value = "runMode"
$ws.Range("B45").Value = "brokerId"
$ws.Range("C45").Value = "agentId"
$ws.Range("D45").Value = "agencyOfficeId"

$ws.Range("A46").Value = "N"
$ws.Range("B46").Value = 20217
$ws.Range("C46").Value = 237
$ws.Range("D46").Value = 8006

$ws.Range("A47").Value = "Y"
$ws.Range("B47").Value = 25997
$ws.Range("C47").Value = 7166
$ws.Range("D47").Value = 8414

# ---------------------------------------------------------------------------
# 3. Update the saved view state of the sheet (active cell A4, no frozen
#    top-left scroll position)
# ---------------------------------------------------------------------------
$ws.Range("A4").Select()
